# Fix issue about 'Configuration of workspaces using system properties'
# Slide 14 ("Configure your workspace using system properties"), shape 2 (body text)
# contains two example system-property names that used a hyphen in
# "repository-collaboration" where an underscore is required:
#   exo.jcr.config.force.workspace.repository-collaboration.container.foo
#   exo.jcr.config.default.workspace.repository-collaboration.container.foo
# Both need "repository-collaboration" -> "repository_collaboration".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Fix #1: bullet "If we have a system property called exo.jcr.config.force.workspace...."
# (select the whole run, including its trailing space, so the run is not split)
$full = $tr.Text
$needle1 = "exo.jcr.config.force.workspace.repository-collaboration.container.foo "
$idx1 = $full.IndexOf($needle1)
if ($idx1 -ge 0) {
    $run1 = $tr.Characters($idx1 + 1, $needle1.Length)
    $run1.Text = "exo.jcr.config.force.workspace.repository_collaboration.container.foo "
}

# --- Fix #2: bullet "If the previous rules don't allow ... exo.jcr.config.default.workspace...."
$full = $tr.Text
$needle2 = "exo.jcr.config.default.workspace.repository-collaboration.container.foo"
$idx2 = $full.IndexOf($needle2)
if ($idx2 -ge 0) {
    $run2 = $tr.Characters($idx2 + 1, $needle2.Length)
    $run2.Text = "exo.jcr.config.default.workspace.repository_collaboration.container.foo"
}
